$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1
$ws.Cells.Item(1,1).Value2 = "Datos actualizados a 27 de Marzo de 2020 a las 12:42"

# Row 4
$ws.Cells.Item(4,2).Value2 = 19243
$ws.Cells.Item(4,3).Value2 = 5044
$ws.Cells.Item(4,4).Value2 = 11787
$ws.Cells.Item(4,5).Value2 = 2412

# Row 7
$ws.Cells.Item(7,1).Value2 = "Valencia/Valencia"
$ws.Cells.Item(7,2).Value2 = 2027
$ws.Cells.Item(7,3).Value2 = 50
$ws.Cells.Item(7,4).Value2 = 1889
$ws.Cells.Item(7,5).Value2 = 88

# Row 8
$ws.Cells.Item(8,1).Value2 = "Bizkaia/Vizcaya"
$ws.Cells.Item(8,2).Value2 = 1850
$ws.Cells.Item(8,3).Value2 = 621
$ws.Cells.Item(8,4).Value2 = 1525
$ws.Cells.Item(8,5).Value2 = 65

# Row 9
$ws.Cells.Item(9,1).Value2 = "Navarra"
$ws.Cells.Item(9,2).Value2 = 1641
$ws.Cells.Item(9,3).Value2 = 70
$ws.Cells.Item(9,4).Value2 = 1513
$ws.Cells.Item(9,5).Value2 = 58

# Row 10
$ws.Cells.Item(10,1).Value2 = "Araba/Alava"
$ws.Cells.Item(10,2).Value2 = 1435
$ws.Cells.Item(10,3).Value2 = 621
$ws.Cells.Item(10,4).Value2 = 1101
$ws.Cells.Item(10,5).Value2 = 94

# Row 11
$ws.Cells.Item(11,1).Value2 = "La Rioja"
$ws.Cells.Item(11,2).Value2 = 1236
$ws.Cells.Item(11,3).Value2 = 62
$ws.Cells.Item(11,4).Value2 = 1119
$ws.Cells.Item(11,5).Value2 = 55

# Row 12
$ws.Cells.Item(12,1).Value2 = "Ciudad Real"
$ws.Cells.Item(12,2).Value2 = 1147
$ws.Cells.Item(12,3).Value2 = 95
$ws.Cells.Item(12,4).Value2 = 1050
$ws.Cells.Item(12,5).Value2 = 89

# Row 13
$ws.Cells.Item(13,1).Value2 = "Alacant/Alicante"
$ws.Cells.Item(13,2).Value2 = 1093
$ws.Cells.Item(13,3).Value2 = 19
$ws.Cells.Item(13,4).Value2 = 983
$ws.Cells.Item(13,5).Value2 = 91

# Row 14
$ws.Cells.Item(14,1).Value2 = "Zaragoza"
$ws.Cells.Item(14,2).Value2 = 1045
$ws.Cells.Item(14,3).Value2 = 68
$ws.Cells.Item(14,4).Value2 = 928
$ws.Cells.Item(14,5).Value2 = 49

# Row 15
$ws.Cells.Item(15,1).Value2 = "Toledo"
$ws.Cells.Item(15,2).Value2 = 965
$ws.Cells.Item(15,3).Value2 = 95
$ws.Cells.Item(15,4).Value2 = 860
$ws.Cells.Item(15,5).Value2 = 78

# Row 16
$ws.Cells.Item(16,1).Value2 = "Aragon"
$ws.Cells.Item(16,2).Value2 = 907
$ws.Cells.Item(16,3).Value2 = 29
$ws.Cells.Item(16,4).Value2 = 838
$ws.Cells.Item(16,5).Value2 = 40

# Row 17
$ws.Cells.Item(17,1).Value2 = "Malaga"
$ws.Cells.Item(17,2).Value2 = 905
$ws.Cells.Item(17,4).Value2 = 811
$ws.Cells.Item(17,5).Value2 = 42

# Row 18
$ws.Cells.Item(18,1).Value2 = "Asturias"
$ws.Cells.Item(18,2).Value2 = 900
$ws.Cells.Item(18,3).Value2 = 52
$ws.Cells.Item(18,4).Value2 = 819
$ws.Cells.Item(18,5).Value2 = 29

# Row 21
$ws.Cells.Item(21,1).Value2 = "Cantabria"
$ws.Cells.Item(21,2).Value2 = 810
$ws.Cells.Item(21,3).Value2 = 19
$ws.Cells.Item(21,4).Value2 = 770
$ws.Cells.Item(21,5).Value2 = 21

# Row 22
$ws.Cells.Item(22,1).Value2 = "Pontevedra"
$ws.Cells.Item(22,2).Value2 = 701
$ws.Cells.Item(22,3).Value2 = 47
$ws.Cells.Item(22,4).Value2 = 684
$ws.Cells.Item(22,5).Value2 = 5

# Row 23
$ws.Cells.Item(23,1).Value2 = "Murcia"
$ws.Cells.Item(23,2).Value2 = 687
$ws.Cells.Item(23,3).Value2 = 12
$ws.Cells.Item(23,4).Value2 = 660
$ws.Cells.Item(23,5).Value2 = 15

# Row 35
$ws.Cells.Item(35,2).Value2 = 412
$ws.Cells.Item(35,4).Value2 = 389
$ws.Cells.Item(35,5).Value2 = 19

# Row 46
$ws.Cells.Item(46,1).Value2 = "Huesca"
$ws.Cells.Item(46,2).Value2 = 150
$ws.Cells.Item(46,3).Value2 = 10
$ws.Cells.Item(46,4).Value2 = 136

# Row 47
$ws.Cells.Item(47,1).Value2 = "Lugo"
$ws.Cells.Item(47,2).Value2 = 145
$ws.Cells.Item(47,3).Value2 = 47
$ws.Cells.Item(47,4).Value2 = 128
$ws.Cells.Item(47,5).Value2 = 4

# Row 48
$ws.Cells.Item(48,1).Value2 = "Teruel"
$ws.Cells.Item(48,2).Value2 = 143
$ws.Cells.Item(48,3).Value2 = 9
$ws.Cells.Item(48,4).Value2 = 129
$ws.Cells.Item(48,5).Value2 = 5

# Row 49
$ws.Cells.Item(49,1).Value2 = "Palencia"
$ws.Cells.Item(49,2).Value2 = 139
$ws.Cells.Item(49,3).Value2 = 14
$ws.Cells.Item(49,4).Value2 = 120
$ws.Cells.Item(49,5).Value2 = 5

# Row 50
$ws.Cells.Item(50,1).Value2 = "Almeria"
$ws.Cells.Item(50,2).Value2 = 134
$ws.Cells.Item(50,3).Value2 = 5
$ws.Cells.Item(50,4).Value2 = 122
$ws.Cells.Item(50,5).Value2 = 7

